$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

# Replace the value in C2 ("ramineni1991" -> "asduhsquac")
$ws.Range("C2").Value = "asduhsquac"

# Move the selection/active cell to C2 (matches the sheetView selection in the diff)
$ws.Activate()
$ws.Range("C2").Select()
